# Task 3 done for task 2
# Fill in the Task 2 (columns H:M) per-philosopher Avg/Max summary rows
# (rows 4-8) and the overall Avg/Max row (row 9), mirroring the Task 1
# summary already present in columns B:G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlEdgeLeft/Top/Bottom/Right + xlLineStyleNone / xlContinuous constants
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlLineStyleNone = -4142
$xlContinuous = 1
$xlThin = 2

# --- Row 4 (Philosopher 1) ---
$ws.Range("H4").Value = 0.6802
$ws.Range("I4").Value = 3.0002
$ws.Range("J4").Value = 0.6573
$ws.Range("K4").Value = 2.0009
$ws.Range("L4").Value = 0.8802
$ws.Range("M4").Value = 5.0003

# --- Row 5 (Philosopher 2) ---
$ws.Range("H5").Value = 0.6001
$ws.Range("I5").Value = 3.0003
$ws.Range("J5").Value = 0.8287
$ws.Range("K5").Value = 3.0009
$ws.Range("L5").Value = 0.5802
$ws.Range("M5").Value = 3.0005

# --- Row 6 (Philosopher 3) ---
$ws.Range("H6").Value = 0.6802
$ws.Range("I6").Value = 3.0006
$ws.Range("J6").Value = 1.0002
$ws.Range("K6").Value = 5.0002
$ws.Range("L6").Value = 1.1403
$ws.Range("M6").Value = 7.0011

# --- Row 7 (Philosopher 4) ---
$ws.Range("H7").Value = 0.4802
$ws.Range("I7").Value = 3.0004
$ws.Range("J7").Value = 0.9431
$ws.Range("K7").Value = 4.0011
$ws.Range("L7").Value = 1.2602
$ws.Range("M7").Value = 7.0007

# --- Row 8 (Philosopher 5) ---
$ws.Range("H8").Value = 0.2803
$ws.Range("I8").Value = 3.0007
$ws.Range("J8").Value = 0.6001
$ws.Range("K8").Value = 3.0005
$ws.Range("L8").Value = 1.1003
$ws.Range("M8").Value = 5.0003

# --- Row 9 totals (Avg. / Max across the five philosophers) ---
$ws.Range("H9").Formula = "=AVERAGE(H4:H8)"
$ws.Range("I9").Formula = "=MAX(I4:I8)"
$ws.Range("J9").Formula = "=AVERAGE(J4:J8)"
$ws.Range("K9").Formula = "=MAX(K4:K8)"
$ws.Range("L9").Formula = "=AVERAGE(L4:L8)"
$ws.Range("M9").Formula = "=MAX(M4:M8)"

# --- Border touch-ups left behind by the original editing session ---
# Each cell below needs exactly ONE border mutation (touching the whole
# Borders collection at once, or a single edge) so the style table is
# reused in place instead of growing with unreferenced intermediate
# styles; multi-area "A1,B1" ranges only affect the first area here, so
# the cell lists are walked individually.

# Row 4 "Max" columns lose their top edge (keep left/right).
foreach ($cell in @("I4","K4","M4")) {
    $ws.Range($cell).Borders.Item($xlEdgeTop).LineStyle = $xlLineStyleNone
}

# Rows 5-7 "Avg" columns lose their border entirely.
foreach ($cell in @("H5","J5","L5","H6","J6","L6","H7","J7","L7")) {
    $ws.Range($cell).Borders.LineStyle = $xlLineStyleNone
}

# Row 7 "Max" columns gain a bottom edge (keep left/right).
foreach ($cell in @("I7","K7","M7")) {
    $b = $ws.Range($cell).Borders.Item($xlEdgeBottom)
    $b.ColorIndex = 1
    $b.LineStyle = $xlContinuous
}

# Row 8 "Avg" columns lose their border entirely.
foreach ($cell in @("H8","J8","L8")) {
    $ws.Range($cell).Borders.LineStyle = $xlLineStyleNone
}

# I8 gains a top edge (keep its existing left/right/bottom -> full box).
$b = $ws.Range("I8").Borders.Item($xlEdgeTop)
$b.ColorIndex = 1
$b.LineStyle = $xlContinuous

# Move the active selection to M14 (where the author's cursor ended up).
$ws.Activate()
$ws.Range("M14").Select()
